$d = $word.ActiveDocument

# The old text (as it reads once all runs in the paragraph are concatenated)
# varies a little bit (leading/trailing spaces) across the four occurrences,
# so we detect the paragraphs generically by looking for the two telltale
# fragments that are common to every variant, then rewrite the whole
# paragraph content with a single plain run containing the new sentence.
$newText = "Datas da campaña de Hercules: 13-22 de xuño, 12-21 de xullo, 10-19 de agosto"

$targets = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Datas da campa*" -and $t -like "*Perseo*" -and $t -like "*decembro*") {
        $targets.Add($p.Range.Start) | Out-Null
    }
}

# Walk the matches back-to-front so earlier offsets stay valid while we
# edit later ones.
for ($k = $targets.Count - 1; $k -ge 0; $k--) {
    $pStart = $targets[$k]
    $p = $d.Range($pStart, $pStart).Paragraphs(1)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    # $pEnd points one past the paragraph/section mark, so the last real
    # character of the paragraph content is at $pEnd - 2 and the mark
    # itself sits at $pEnd - 1.
    $textEnd = $pEnd - 1

    # A couple of the paragraphs end with a dangling w:proofErr element
    # that sits right up against the paragraph mark, after the very last
    # run. A plain Range.Delete() over the text can't reach/clear such a
    # marker when nothing follows it inside the paragraph, so we first
    # insert a throwaway character right after the existing content -
    # this creates a run after the proofErr marker. Including that
    # throwaway character in the subsequent delete then correctly clears
    # any stray proofErr markers along with all the old runs.
    $dummy = $d.Range($textEnd, $textEnd)
    $dummy.InsertAfter("X")

    $full = $d.Range($pStart, $textEnd + 1)
    $full.Delete()

    $ins = $d.Range($pStart, $pStart)
    $ins.InsertAfter($newText)
}
